# logBook.xlsx - add new log entry (row 16) for "FCN_starter pipeline" work,
# matching formatting of the preceding entry (row 15), and update the
# active selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row of data at row 16, taking on the same number
#     formats / alignment as row 15 immediately above it (date, time,
#     time, duration formula, centered Sno, wrapped Description). ---
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 44728
$ws.Range("C16").Value = 0.125
$ws.Range("D16").Value = 0.26041666666666669
$ws.Range("E16").Formula = "=D16-C16"
$ws.Range("F16").Value = "Code"
$ws.Range("G16").Value = "Pipeline for FCN starter notebook completed;`nNOTE : need to change labels, resize image, label, add class weights in nn.CE, restructure notebook last section"

$ws.Range("G16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 45

# --- Update the sheet view: scroll down a bit and move the selection
#     to the newly entered description cell. ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("G16").Select()
